$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-04 Sunday" "2024-08-05 Monday"

Replace-Text "40×98=3920" "19×41=779"
Replace-Text "35×70=2450" "15×56=840"
Replace-Text "18×30=540" "18×20=360"
Replace-Text "58×18=1044" "53×32=1696"
Replace-Text "29×30=870" "43×76=3268"

Replace-Text "32×24=768" "21×89=1869"
Replace-Text "77×90=6930" "73×49=3577"
Replace-Text "65×23=1495" "92×32=2944"
Replace-Text "32×29=928" "89×77=6853"
Replace-Text "87×53=4611" "68×14=952"

Replace-Text "78×14=1092" "13×68=884"
Replace-Text "46×28=1288" "42×78=3276"
Replace-Text "70×14=980" "25×71=1775"
Replace-Text "58×92=5336" "38×73=2774"
Replace-Text "86×14=1204" "71×61=4331"

Replace-Text "82×13=1066" "74×68=5032"
Replace-Text "70×83=5810" "39×75=2925"
Replace-Text "54×70=3780" "63×40=2520"
Replace-Text "77×60=4620" "17×36=612"
Replace-Text "85×76=6460" "97×56=5432"

Replace-Text "14×82=1148" "39×42=1638"
Replace-Text "39×89=3471" "82×63=5166"
Replace-Text "53×42=2226" "74×78=5772"
Replace-Text "92×47=4324" "49×18=882"
Replace-Text "15×82=1230" "75×58=4350"
